# Update workbook view zoom/window height slightly
$wb = $excel.ActiveWorkbook

# --- Workbook-level view tweak ---
$wb.Windows.Item(1).Height = 14780

# --- Sheet1: populate new av_fruit_time (H) and shade_tol (I) data, plus a few corrected values ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Cells.Item(2, 8).Value = 9.5
$ws1.Cells.Item(3, 8).Value = 9.5
$ws1.Cells.Item(4, 8).Value = 9.5
$ws1.Cells.Item(5, 8).Value = 10
$ws1.Cells.Item(8, 8).Value = 5.5
$ws1.Cells.Item(8, 9).Value = "very_intolerant"
$ws1.Cells.Item(9, 9).Value = "intolerant"
$ws1.Cells.Item(10, 8).Value = 5.5
$ws1.Cells.Item(10, 9).Value = "intolerant"
$ws1.Cells.Item(11, 8).Value = 5.5
$ws1.Cells.Item(11, 9).Value = "intolerant"
$ws1.Cells.Item(12, 8).Value = 5.5
$ws1.Cells.Item(12, 9).Value = "very_intolerant"
$ws1.Cells.Item(13, 8).Value = 5.5
$ws1.Cells.Item(13, 9).Value = "intolerant"
$ws1.Cells.Item(14, 8).Value = 5.5
$ws1.Cells.Item(14, 9).Value = "very_intolerant"
$ws1.Cells.Item(15, 8).Value = 5.5
$ws1.Cells.Item(15, 9).Value = "very_intolerant"
$ws1.Cells.Item(16, 8).Value = 5
$ws1.Cells.Item(16, 9).Value = "intolerant"
$ws1.Cells.Item(17, 8).Value = 5.5
$ws1.Cells.Item(17, 9).Value = "very_intolerant"
$ws1.Cells.Item(18, 8).Value = 6
$ws1.Cells.Item(18, 9).Value = "very_intolerant"
$ws1.Cells.Item(19, 8).Value = 6.5
$ws1.Cells.Item(19, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(20, 8).Value = 9.5
$ws1.Cells.Item(20, 9).Value = "intolerant"
$ws1.Cells.Item(21, 8).Value = 9.5
$ws1.Cells.Item(21, 9).Value = "intolerant"
$ws1.Cells.Item(22, 8).Value = 10
$ws1.Cells.Item(22, 9).Value = "very_intolerant"
$ws1.Cells.Item(24, 8).Value = 10
$ws1.Cells.Item(24, 9).Value = "intolerant"
$ws1.Cells.Item(25, 8).Value = 8.5
$ws1.Cells.Item(25, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(26, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(27, 8).Value = 8.5
$ws1.Cells.Item(27, 9).Value = "intolerant"
$ws1.Cells.Item(28, 8).Value = 7.5
$ws1.Cells.Item(28, 9).Value = "intolerant"
$ws1.Cells.Item(29, 8).Value = 8.5
$ws1.Cells.Item(29, 9).Value = "tolerant"
$ws1.Cells.Item(30, 8).Value = 7.5
$ws1.Cells.Item(30, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(31, 8).Value = 8.5
$ws1.Cells.Item(31, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(33, 8).Value = 5
$ws1.Cells.Item(33, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(34, 8).Value = 5
$ws1.Cells.Item(34, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(35, 8).Value = 5
$ws1.Cells.Item(35, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(36, 8).Value = 4
$ws1.Cells.Item(36, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(37, 8).Value = 9.5
$ws1.Cells.Item(37, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(38, 8).Value = 9
$ws1.Cells.Item(38, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(39, 8).Value = 9.5
$ws1.Cells.Item(39, 9).Value = "intolerant"
$ws1.Cells.Item(40, 8).Value = 6.5
$ws1.Cells.Item(40, 9).Value = "intolerant"
$ws1.Cells.Item(41, 8).Value = 7
$ws1.Cells.Item(41, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(42, 8).Value = 9.5
$ws1.Cells.Item(42, 9).Value = "intolerant"
$ws1.Cells.Item(43, 8).Value = 9.5
$ws1.Cells.Item(43, 9).Value = "very_tolerant"
$ws1.Cells.Item(44, 8).Value = 9.5
$ws1.Cells.Item(44, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(45, 8).Value = 9.5
$ws1.Cells.Item(45, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(46, 8).Value = 9.5
$ws1.Cells.Item(46, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(47, 8).Value = 9.5
$ws1.Cells.Item(47, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(48, 8).Value = 19
$ws1.Cells.Item(49, 8).Value = 19
$ws1.Cells.Item(49, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(50, 8).Value = 19
$ws1.Cells.Item(50, 9).Value = "intolerant"
$ws1.Cells.Item(51, 8).Value = 19
$ws1.Cells.Item(51, 9).Value = "intolerant"
$ws1.Cells.Item(52, 8).Value = 19
$ws1.Cells.Item(52, 9).Value = "intolerant"
$ws1.Cells.Item(53, 8).Value = 19
$ws1.Cells.Item(53, 9).Value = "intolerant"
$ws1.Cells.Item(54, 8).Value = 19
$ws1.Cells.Item(54, 9).Value = "intolerant"
$ws1.Cells.Item(55, 8).Value = 19
$ws1.Cells.Item(55, 9).Value = "intolerant"
$ws1.Cells.Item(57, 8).Value = 9.5
$ws1.Cells.Item(57, 9).Value = "intolerant"
$ws1.Cells.Item(58, 8).Value = 8
$ws1.Cells.Item(58, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(59, 8).Value = 8.5
$ws1.Cells.Item(59, 9).Value = "very_intolerant"
$ws1.Cells.Item(60, 8).Value = 8.5
$ws1.Cells.Item(60, 9).Value = "very_intolerant"
$ws1.Cells.Item(61, 8).Value = 5.5
$ws1.Cells.Item(61, 9).Value = "very_intolerant"
$ws1.Cells.Item(62, 8).Value = 6
$ws1.Cells.Item(62, 9).Value = "very_tolerant"
$ws1.Cells.Item(63, 8).Value = 9
$ws1.Cells.Item(63, 9).Value = "tolerant"
$ws1.Cells.Item(64, 8).Value = 10
$ws1.Cells.Item(64, 9).Value = "tolerant"
$ws1.Cells.Item(65, 8).Value = 10
$ws1.Cells.Item(65, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(66, 8).Value = 10
$ws1.Cells.Item(66, 9).Value = "very_tolerant"
$ws1.Cells.Item(67, 8).Value = 10
$ws1.Cells.Item(67, 9).Value = "tolerant"
$ws1.Cells.Item(68, 8).Value = 10
$ws1.Cells.Item(68, 9).Value = "intolerant"
$ws1.Cells.Item(69, 8).Value = 10
$ws1.Cells.Item(69, 9).Value = "intolerant"
$ws1.Cells.Item(70, 8).Value = 10
$ws1.Cells.Item(70, 9).Value = "tolerant"
$ws1.Cells.Item(71, 8).Value = 10
$ws1.Cells.Item(71, 9).Value = "tolerant"
$ws1.Cells.Item(72, 8).Value = 10
$ws1.Cells.Item(72, 9).Value = "intolerant"
$ws1.Cells.Item(75, 8).Value = 10
$ws1.Cells.Item(75, 9).Value = "very_tolerant"
$ws1.Cells.Item(76, 8).Value = 9.5
$ws1.Cells.Item(76, 9).Value = "very_tolerant"
$ws1.Cells.Item(77, 4).Value = "wind"
$ws1.Cells.Item(77, 8).Value = 5.5
$ws1.Cells.Item(77, 9).Value = "tolerant"
$ws1.Cells.Item(78, 4).Value = "wind"
$ws1.Cells.Item(78, 8).Value = 5
$ws1.Cells.Item(78, 9).Value = "intolerant"
$ws1.Cells.Item(81, 4).Value = "wind"
$ws1.Cells.Item(81, 8).Value = 9.5
$ws1.Cells.Item(81, 9).Value = "very_tolerant"
$ws1.Cells.Item(82, 4).Value = "wind"
$ws1.Cells.Item(82, 8).Value = 9.5
$ws1.Cells.Item(82, 9).Value = "intolerant"
$ws1.Cells.Item(83, 8).Value = 10
$ws1.Cells.Item(83, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(84, 7).Value = "perfect"
$ws1.Cells.Item(84, 8).Value = 10
$ws1.Cells.Item(84, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(86, 8).Value = 10
$ws1.Cells.Item(86, 9).Value = "very_tolerant"
$ws1.Cells.Item(87, 8).Value = 10
$ws1.Cells.Item(87, 9).Value = "tolerant"
$ws1.Cells.Item(88, 8).Value = 8.5
$ws1.Cells.Item(88, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(89, 8).Value = 8.5
$ws1.Cells.Item(89, 9).Value = "moderately_tolerant"
$ws1.Cells.Item(90, 8).Value = 9
$ws1.Cells.Item(90, 9).Value = "intolerant"
$ws1.Cells.Item(91, 8).Value = 9
$ws1.Cells.Item(91, 9).Value = "intolerant"
$ws1.Cells.Item(92, 8).Value = 9.5
$ws1.Cells.Item(92, 9).Value = "tolerant"
$ws1.Cells.Item(93, 8).Value = 9
$ws1.Cells.Item(93, 9).Value = "intolerant"
$ws1.Cells.Item(94, 8).Value = 9
$ws1.Cells.Item(94, 9).Value = "intolerant"

# --- Sheet1 view: scroll position + active selection ---
$ws1.Application.ActiveWindow.ScrollRow = 70
$ws1.Range("G85").Select()

# --- Sheet2: add new Metadata legend columns for av.fruit time (D) and shade_tol (E) ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Cells.Item(2, 4).Value = "av.fruit time"
$ws2.Cells.Item(2, 4).Font.Bold = $true
$ws2.Cells.Item(3, 4).Value = "1-12=Jan-Dec"
$ws2.Cells.Item(4, 4).Value = "autmn=9.5"
$ws2.Cells.Item(5, 4).Value = "late autumn=10"
$ws2.Cells.Item(6, 4).Value = "early autmn=9"
$ws2.Cells.Item(2, 5).Value = "shade_tol"
$ws2.Cells.Item(2, 5).Font.Bold = $true
$ws2.Cells.Item(3, 5).Value = "very intolerant"
$ws2.Cells.Item(4, 5).Value = "intolerant"
$ws2.Cells.Item(5, 5).Value = "moderately tolerant"
$ws2.Cells.Item(6, 5).Value = "tolerant"
$ws2.Cells.Item(7, 5).Value = "very tolerant"

# --- Sheet2 view: active selection ---
$ws2.Range("E8").Select()
